$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CO2")

# Rename the "Storage" column header to "Battery"
$ws.Range("G2").Value = "Battery"

# Rename the matching defined name so it stays in sync with the new label
$wb.Names.Item("CO2_gen_Storage").Delete()
$wb.Names.Add("CO2_gen_Battery", "='CO2'!`$G`$5")

# Match the author's last-touched cell/selection
$ws.Range("G5").Select()
